$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 354.5
    3  = 342.3
    4  = 340.7
    5  = 300.2
    6  = 305.1
    7  = 295.3
    8  = 389.9
    9  = 419.3
    10 = 326.2
    11 = 343.1
    12 = 336.7
    13 = 311.1
    14 = 402.9
    15 = 417.8
    16 = 396.3
    17 = 424.2
    18 = 318.7
    19 = 345.9
    20 = 435.1
    21 = 417.2
    22 = 437.6
    23 = 466.3
    24 = 434.6
    25 = 426.9
    26 = 521.5
    27 = 518.5
    28 = 475.9
    29 = 449
    30 = 512.9
    31 = 465.4
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
